# Apply the "Updated Basic Weapons and Enemies again" edit.
#
# This adds "average"/"fast"/"slow" movement-speed call-outs to each enemy
# bullet, fixes "but" -> "and" for the Staff Ogre, and appends a note about
# the Sword possibly reflecting/destroying projectiles (merging the Sword
# bullet with the trailing bookmark-only paragraph, then leaving a fresh
# blank paragraph behind, like Word does when you backspace a paragraph
# break and then press Enter again later).

$d = $word.ActiveDocument

function Insert-BeforeMatch {
    param(
        [string]$matchText,
        [string]$insertText
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($matchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $matchText"
    }
    $rng.Collapse(0)
    $rng.MoveEnd(1, -1)
    $rng.InsertBefore($insertText)
}

# 1. Standard Skeleton
Insert-BeforeMatch "Standard Skeleton (Moves towards player, uses standard weapon, low health)" " and average movement speed"

# 2. Mage
Insert-BeforeMatch "Mage (Moves towards player, uses Burst-Fire weapon, average health)" " and average movement speed"

# 3. Ghost
Insert-BeforeMatch "Ghost (Moves away from player, uses Rapid-Fire weapon, low health)" " and average movement speed"

# 4. Staff Ogre: "but" -> "and", then append " speed" before the closing paren
$d.Content.Find.Execute("high health but slow movement)", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "high health and slow movement)", 2) | Out-Null
Insert-BeforeMatch "Staff Ogre (Moves towards player, uses Staff, high health and slow movement)" " speed"

# 5. Pyromaniac
Insert-BeforeMatch "Pyromaniac (Moves towards player quickly, uses Flamethrower, low health and fast movement)" " speed"

# 6. Necromancer
Insert-BeforeMatch "Necromancer (Moves away from player slowly, spawns Standard Skeletons and shoots a ring of projectiles around itself occasionally, High health and slow movement)" " speed"

# 7. Knight: paragraph ends right after "movement" (no closing paren in the source)
$knightPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq "Knight (Moves towards player, uses Sword (does not fire projectiles), High health and average movement") {
        $knightPara = $p
        break
    }
}
if ($null -eq $knightPara) {
    throw "Could not find Knight paragraph"
}
$kr = $knightPara.Range.Duplicate
$kr.Collapse(0)
$kr.MoveEnd(1, -1)
$kr.InsertAfter(" speed)")

# 8. Sword (weapon list): append the projectile-reflect note, then merge the
#    trailing bookmark-only paragraph into it while keeping the Sword
#    paragraph's own list formatting, and leave a new blank paragraph behind.
Insert-BeforeMatch "Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage)" ", (maybe reflects/destroys projectiles it hits?)"

$swordPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage, (maybe reflects/destroys projectiles it hits?)")) {
        $swordPara = $p
        break
    }
}
if ($null -eq $swordPara) {
    throw "Could not find Sword paragraph"
}

# Re-home the hidden _GoBack bookmark onto the spot right before the closing
# ")" we're about to add (Find-derived ranges collapse to a reliable zero
# width point; Bookmarks.Add needs that to avoid snapping to the paragraph
# start).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmRange = $d.Content
$bmRange.Find.Execute("hits?)", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bmRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# Merge the now-empty bookmark paragraph into the Sword paragraph by deleting
# the paragraph mark between them (the survivor keeps the Sword paragraph's
# own list-paragraph formatting because the bookmark paragraph carried none).
$swordPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage, (maybe reflects/destroys projectiles it hits?)")) {
        $swordPara = $p
        break
    }
}
$sr = $swordPara.Range
$markRange = $d.Range($sr.End - 1, $sr.End)
$markRange.Delete()

# Add the closing paren after the bookmark.
$swordPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage, (maybe reflects/destroys projectiles it hits?)")) {
        $swordPara = $p
        break
    }
}
$sr2 = $swordPara.Range.Duplicate
$sr2.Collapse(0)
$sr2.MoveEnd(1, -1)
$sr2.InsertAfter(")")

# Leave a fresh blank paragraph behind, like the diff shows.
$swordPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Sword (Melee weapon, hits enemies directly in front of player, fast attack speed, high damage, (maybe reflects/destroys projectiles it hits?)")) {
        $swordPara = $p
        break
    }
}
$endRng = $swordPara.Range.Duplicate
$endRng.Collapse(0)
$endRng.InsertAfter([char]13)

Write-Output "done"
